# Applies the commit's change: on the "DeviceList" sheet, delete columns
# B, E, F and G (the four devices that are no longer part of this batch),
# which shifts the remaining device columns (old C,D,H,I) left into
# B,C,D,E. Also update the recorded absolute path in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

# Delete columns right-to-left so earlier column letters stay valid
# while we work through the list.
$ws.Range("G:G").Delete() | Out-Null
$ws.Range("F:F").Delete() | Out-Null
$ws.Range("E:E").Delete() | Out-Null
$ws.Range("B:B").Delete() | Out-Null

$ws.Range("E17").Select() | Out-Null
